# Generate Report for Handback
# Updates row 7 (b8f7c47a-e173-4672-a45d-e202b23712cb) on both the
# "zh-cn" and "de-de" sheets: the handback now arrived, so the
# "Latest Target File", "Latest Handback File" and "Latest Handback
# DateTime" columns get filled in, and an Error Detail is recorded
# because the handback was produced against a stale commit.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/140c9718e668f5fae8f99bcb34912654b8ea9beb/e2e/b8f7c47a-e173-4672-a45d-e202b23712cb.md"
$targetDisplay = "b8f7c47a-e173-4672-a45d-e202b23712cb.md"

function Update-Row7 {
    param(
        $ws,
        [string]$handbackFile,
        [string]$handbackDateTime
    )

    $errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/21e84dc15ac431912dbe91e6952c9c3f373a6aeb/e2e/b8f7c47a-e173-4672-a45d-e202b23712cb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/140c9718e668f5fae8f99bcb34912654b8ea9beb/e2e/b8f7c47a-e173-4672-a45d-e202b23712cb.md."

    # Latest Target File (I7): becomes a hyperlink to the source md file,
    # matching the style/behaviour already used by column A and the other
    # "Latest Target File" cells (I2..I5).
    $ws.Hyperlinks.Add($ws.Range("I7"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetDisplay) | Out-Null

    # Make sure the cell picks up the same visual "HyperLink" look (blue,
    # underlined) used everywhere else in this workbook.
    $ws.Range("I7").Style = "HyperLink"
    $ws.Range("I7").Font.Color = 15570276
    $ws.Range("I7").Font.Underline = 2

    # Latest Handback File (J7)
    $ws.Range("J7").Value = $handbackFile

    # Latest Handback DateTime (K7)
    $ws.Range("K7").Value = $handbackDateTime

    # Error Detail (P7)
    $ws.Range("P7").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-Row7 $wsZhCn "b8f7c47a-e173-4672-a45d-e202b23712cb.2bfbe02784bd9bbfe0f5988ab6057408be837234.zh-cn.xlf" "2016-09-05 21:04:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-Row7 $wsDeDe "b8f7c47a-e173-4672-a45d-e202b23712cb.2bfbe02784bd9bbfe0f5988ab6057408be837234.de-de.xlf" "2016-09-05 21:04:24"
